# Weekly refresh: insert the latest week's two price observations
# (Primera / Segunda calidad) at the top of the Zapallo italiano data
# block, pushing the existing 131:191 rows down to 133:193.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("131:132").Insert()

# New row 131 - Primera
$ws.Range("A131").Value = 9
$ws.Range("B131").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C131").Value = "Metropolitana"
$ws.Range("D131").Value = 44452
$ws.Range("E131").Value = 13
$ws.Range("F131").Value = 100112032
$ws.Range("G131").Value = "Zapallo italiano"
$ws.Range("H131").Value = "Sin especificar"
$ws.Range("I131").Value = "Primera"
$ws.Range("J131").Value = 79
$ws.Range("K131").Value = 17000
$ws.Range("L131").Value = 18000
$ws.Range("M131").Value = 17494
$ws.Range("N131").Value = "$/caja 60 unidades"
$ws.Range("O131").Value = "Región de Arica y Parinacota"
$ws.Range("P131").Value = 292
$ws.Range("Q131").Value = 60
$ws.Range("R131").Value = "Hortaliza"

# New row 132 - Segunda
$ws.Range("A132").Value = 9
$ws.Range("B132").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C132").Value = "Metropolitana"
$ws.Range("D132").Value = 44452
$ws.Range("E132").Value = 13
$ws.Range("F132").Value = 100112032
$ws.Range("G132").Value = "Zapallo italiano"
$ws.Range("H132").Value = "Sin especificar"
$ws.Range("I132").Value = "Segunda"
$ws.Range("J132").Value = 43
$ws.Range("K132").Value = 16000
$ws.Range("L132").Value = 16000
$ws.Range("M132").Value = 16000
$ws.Range("N132").Value = "$/caja 100 unidades"
$ws.Range("O132").Value = "Región de Arica y Parinacota"
$ws.Range("P132").Value = 160
$ws.Range("Q132").Value = 100
$ws.Range("R132").Value = "Hortaliza"
